$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("A1").Value = "Keyword"
$ws.Range("B1").Value = "Input1"
$ws.Range("C1").Value = "Input2"
$ws.Range("D1").Value = "Input3"

# --- Row 2: test case id ---
$ws.Range("A2").Value = "tc_id"
$ws.Range("D2").Value = "TS001"

# --- Row 3: test case description ---
$ws.Range("A3").Value = "tc_desc"
$ws.Range("D3").Value = "To verify if the below:
1. All checkboxes able to check."
$ws.Rows.Item(3).RowHeight = 28.8

# --- Row 4: first step (open browser) ---
$ws.Range("A4").Value = "step"
$ws.Range("B4").Value = "Open the browser"
$ws.Range("C4").Value = "The browser opens successfully"

# --- Row 5: open_browser keyword ---
$ws.Range("A5").Value = "open_browser"
$ws.Range("B5").Value = "Chrome"
$ws.Range("D5").Value = "Chrome"

# --- Row 6: enter_url keyword, now pointed at jqueryui.com ---
$ws.Range("A6").Value = "enter_url"
$ws.Range("D6").Value = "https://jqueryui.com/"

# --- Remove old leftover rows 7 through 42 (stale steps + stray formatted cells) ---
$ws.Range("A7:D42").EntireRow.Delete()

# --- Row 7: open checkbox radio page ---
$ws.Range("A7").Value = "step"
$ws.Range("B7").Value = "Open the checkbox radio page"
$ws.Range("C7").Value = "The checkbox radio page opens successfully"

# --- Row 8: click checkboxradio link ---
$ws.Range("A8").Value = "click"
$ws.Range("B8").Value = "Checkboxradio link"
$ws.Range("C8").Value = "checkbox_page_css"

# --- Row 9: step - click on no icons link ---
$ws.Range("A9").Value = "step"
$ws.Range("B9").Value = "Click on no icons link"
$ws.Range("C9").Value = "The no icons page opens successfully"

# --- Row 10: click checkboxnoicons link ---
$ws.Range("A10").Value = "click"
$ws.Range("B10").Value = "checkboxnoicons link"
$ws.Range("C10").Value = "checkbox_no_icons_link_css"

# --- Row 11: step - click on all radio buttons ---
$ws.Range("A11").Value = "step"
$ws.Range("B11").Value = "Click on all radio buttons anc check"
$ws.Range("C11").Value = "Should be able to click on each radio button and that button is selected"

# --- Row 12: new switch_to_iframe keyword ---
$ws.Range("A12").Value = "switch_to_iframe"
$ws.Range("B12").Value = "iframe"
$ws.Range("C12").Value = "checkbox_iframe_css"

# --- Row 13: click 1radio ---
$ws.Range("A13").Value = "click"
$ws.Range("B13").Value = "1radio"
$ws.Range("C13").Value = "ny_radio_css"

# --- Row 14: click 2radio ---
$ws.Range("A14").Value = "click"
$ws.Range("B14").Value = "2radio"
$ws.Range("C14").Value = "paris_radio_css"

# --- Row 15: click 3radio ---
$ws.Range("A15").Value = "click"
$ws.Range("B15").Value = "3radio"
$ws.Range("C15").Value = "london_radio_css"

# --- Row 16: step - click on all check boxes ---
$ws.Range("A16").Value = "step"
$ws.Range("B16").Value = "Click on all check boxes anc check"
$ws.Range("C16").Value = "Should be able to click on all check boxes and all should be checked."

# --- Row 17: click 1check ---
$ws.Range("A17").Value = "click"
$ws.Range("B17").Value = "1check"
$ws.Range("C17").Value = "2star_check_css"

# --- Row 18: click 2check ---
$ws.Range("A18").Value = "click"
$ws.Range("B18").Value = "2check"
$ws.Range("C18").Value = "3star_check_css"

# --- Row 19: click 3check ---
$ws.Range("A19").Value = "click"
$ws.Range("B19").Value = "3check"
$ws.Range("C19").Value = "4star_check_css"

# --- Row 20: click 4check ---
$ws.Range("A20").Value = "click"
$ws.Range("B20").Value = "4check"
$ws.Range("C20").Value = "5star_check_css"

# --- Stray formatted (but empty) cells left over further down column D ---
$ws.Range("D24").WrapText = $true

$ws.Range("D26").WrapText = $true

$ws.Range("D29").NumberFormat = "@"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").WrapText = $true

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").WrapText = $true

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").WrapText = $true

$ws.Range("D37").NumberFormat = "@"

$ws.Range("D38").NumberFormat = "@"

$ws.Range("D39").NumberFormat = "@"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").WrapText = $true

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").WrapText = $true

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").WrapText = $true

$ws.Range("D47").NumberFormat = "@"

# --- Selection matches the author's final cursor position ---
$ws.Range("C17").Select()
